# Site updated: 2021-04-17 12:08:21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (row 19 and row 21, row 20 left empty) ---
$ws.Cells.Item(19, 1).Value = 139
$ws.Cells.Item(19, 2).Value = "拆分词句"
$ws.Cells.Item(19, 4).Value = "||"
$ws.Cells.Item(19, 5).Value = "可以用动态规划，也可以用回溯，比较经典"

$ws.Cells.Item(21, 1).Value = 141
$ws.Cells.Item(21, 2).Value = "单链表中的环"
$ws.Cells.Item(21, 4).Value = "|||"
$ws.Cells.Item(21, 5).Value = "快慢指针，用了比较巧妙的方法"

# --- View state: scroll + selection mirror the saved workbook state ---
$ws.Range("E21").Select()
$excel.ActiveWindow.ScrollRow = 6

# --- Window size tweak recorded in workbook.xml ---
$excel.ActiveWindow.Height = 13200
